$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.019
$ws.Range("A4").Value = -21.774
$ws.Range("D4").Value = -8.134
$ws.Range("E4").Value = 13.019
$ws.Range("D5").Value = -8.616999999999999
$ws.Range("A6").Value = -21.14
$ws.Range("D6").Value = -8.529
$ws.Range("A7").Value = -21.018
$ws.Range("A8").Value = -21.018
$ws.Range("D8").Value = -8.400000000000002
$ws.Range("E9").Value = 12.961
$ws.Range("E11").Value = 12.774
$ws.Range("E14").Value = 13.06
$ws.Range("A16").Value = -20.719
$ws.Range("D16").Value = -8.427
$ws.Range("E18").Value = 12.596
$ws.Range("A20").Value = -21.86
$ws.Range("A21").Value = -20.921
$ws.Range("D22").Value = -8.134
$ws.Range("E25").Value = 12.791
